# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job
# profit sheets, matching the latest Universalis price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1001137.2
$ws.Range("I33").Value = 1251346.5
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 1251346.5
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = -1251117.5
$ws.Range("N33").Value = -758

$ws.Range("H43").Value = 1982.1111
$ws.Range("I43").Value = 2250
$ws.Range("J43").Value = 1767.8
$ws.Range("K43").Value = 2250
$ws.Range("L43").Value = 1767.8
$ws.Range("M43").Value = -2181
$ws.Range("N43").Value = -1905.8

$ws.Range("N58").ClearContents()
$ws.Range("H58").Value = 207.5
$ws.Range("I58").Value = 207.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 622.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -472.5

$ws.Range("M117").ClearContents()
$ws.Range("H117").Value = 40000
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 40000
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 40000
$ws.Range("N117").Value = -49178

$ws.Range("H132").Value = 12558.571
$ws.Range("I132").Value = 8902.777
$ws.Range("J132").Value = 34493.332
$ws.Range("K132").Value = 26708.331
$ws.Range("L132").Value = 103479.996
$ws.Range("M132").Value = -24178.331
$ws.Range("N132").Value = -108539.996

$ws.Range("H138").Value = 2153.875
$ws.Range("I138").Value = 1082.1111
$ws.Range("J138").Value = 3151.724
$ws.Range("K138").Value = 3246.3333
$ws.Range("L138").Value = 9455.172
$ws.Range("M138").Value = 1893.6667
$ws.Range("N138").Value = -19735.172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2435.55
$ws.Range("I61").Value = 2726.75
$ws.Range("J61").Value = 1998.75
$ws.Range("K61").Value = 2726.75
$ws.Range("L61").Value = 1998.75
$ws.Range("M61").Value = -2514.75
$ws.Range("N61").Value = -2422.75

$ws.Range("N97").ClearContents()
$ws.Range("H97").Value = 333.75
$ws.Range("I97").Value = 333.75
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 333.75
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 162.25

$ws.Range("H136").Value = 2435.55
$ws.Range("I136").Value = 2726.75
$ws.Range("J136").Value = 1998.75
$ws.Range("K136").Value = 8180.25
$ws.Range("L136").Value = 5996.25
$ws.Range("M136").Value = -5630.25
$ws.Range("N136").Value = -11096.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 23000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 23000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 23000
$ws.Range("N2").Value = -23226

$ws.Range("H107").Value = 866.6667
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 920
$ws.Range("N107").Value = -4640

$ws.Range("H134").Value = 10278.134
$ws.Range("I134").Value = 6913.143
$ws.Range("J134").Value = 18129.777
$ws.Range("K134").Value = 20739.429
$ws.Range("L134").Value = 54389.33099999999
$ws.Range("M134").Value = -18204.429
$ws.Range("N134").Value = -59459.33099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 5065.8096
$ws.Range("I7").Value = 7163.357
$ws.Range("J7").Value = 870.7143
$ws.Range("K7").Value = 7163.357
$ws.Range("L7").Value = 870.7143
$ws.Range("M7").Value = -7050.357
$ws.Range("N7").Value = -1096.7143

$ws.Range("H62").Value = 13745.179
$ws.Range("I62").Value = 2509.75
$ws.Range("J62").Value = 41833.75
$ws.Range("K62").Value = 2509.75
$ws.Range("L62").Value = 41833.75
$ws.Range("M62").Value = -1885.75
$ws.Range("N62").Value = -43081.75

$ws.Range("H65").Value = 13745.179
$ws.Range("I65").Value = 2509.75
$ws.Range("J65").Value = 41833.75
$ws.Range("K65").Value = 12548.75
$ws.Range("L65").Value = 209168.75
$ws.Range("M65").Value = -9428.75
$ws.Range("N65").Value = -215408.75

$ws.Range("H107").Value = 1194.0526
$ws.Range("I107").Value = 633.7
$ws.Range("J107").Value = 1816.6666
$ws.Range("K107").Value = 633.7
$ws.Range("L107").Value = 1816.6666
$ws.Range("M107").Value = 1286.3
$ws.Range("N107").Value = -5656.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 969.475
$ws.Range("I107").Value = 506
$ws.Range("J107").Value = 981.35895
$ws.Range("K107").Value = 1518
$ws.Range("L107").Value = 2944.07685
$ws.Range("M107").Value = 402
$ws.Range("N107").Value = -6784.07685

$ws.Range("H131").Value = 587452.9399999999
$ws.Range("I131").Value = 549.1667
$ws.Range("J131").Value = 1001737.94
$ws.Range("K131").Value = 1647.5001
$ws.Range("L131").Value = 3005213.82
$ws.Range("M131").Value = 3392.4999
$ws.Range("N131").Value = -3015293.82

$ws.Range("H140").Value = 1582.2727
$ws.Range("I140").Value = 1165.2941
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 3495.8823
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = 1684.1177
$ws.Range("N140").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M40").ClearContents()
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5302

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 795.8889
$ws.Range("I22").Value = 793.8333
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 793.8333
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -498.8333
$ws.Range("N22").Value = -1390

$ws.Range("H27").Value = 795.8889
$ws.Range("I27").Value = 793.8333
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 793.8333
$ws.Range("L27").Value = 800
$ws.Range("M27").Value = -686.8333
$ws.Range("N27").Value = -1014

$ws.Range("H46").Value = 589371.7
$ws.Range("I46").Value = 758.8889
$ws.Range("J46").Value = 1251561.1
$ws.Range("K46").Value = 758.8889
$ws.Range("L46").Value = 1251561.1
$ws.Range("M46").Value = -570.8889
$ws.Range("N46").Value = -1251937.1

$ws.Range("H68").Value = 3414.8572
$ws.Range("I68").Value = 3301.3333
$ws.Range("J68").Value = 3500
$ws.Range("K68").Value = 3301.3333
$ws.Range("L68").Value = 3500
$ws.Range("M68").Value = -2552.3333
$ws.Range("N68").Value = -4998

$ws.Range("H71").Value = 3414.8572
$ws.Range("I71").Value = 3301.3333
$ws.Range("J71").Value = 3500
$ws.Range("K71").Value = 16506.6665
$ws.Range("L71").Value = 17500
$ws.Range("M71").Value = -12762.6665
$ws.Range("N71").Value = -24988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 10000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 10000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 10000
$ws.Range("N13").Value = -10280

$ws.Range("H117").Value = 23266.666
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 23266.666
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 23266.666
$ws.Range("N117").Value = -32444.666
